$d = $word.ActiveDocument

# ---------------------------------------------------------------------
# 1. Insert a new header block at the very top of the document:
#      "Unknown Name"                                   (Title style)
#      "Email: unknown@example.com | Phone: Unknown Phone" (Normal)
#      (empty paragraph)                                (Normal)
# ---------------------------------------------------------------------
$styleNormal = $d.Styles("Normal")

$firstPar = $d.Paragraphs(1)
$firstPar.Range.InsertParagraphBefore()
$titlePar = $d.Paragraphs(1)
$titlePar.Range.Text = "Unknown Name"
$titlePar.Style = "Title"

$secondPar = $d.Paragraphs(2)
$secondPar.Range.InsertParagraphBefore()
$contactPar = $d.Paragraphs(2)
$contactPar.Range.Text = "Email: unknown@example.com | Phone: Unknown Phone"
$contactPar.Style = $styleNormal

$thirdPar = $d.Paragraphs(3)
$thirdPar.Range.InsertParagraphBefore()
$blankPar = $d.Paragraphs(3)
$blankPar.Style = $styleNormal

# ---------------------------------------------------------------------
# 2. Rewrite the Professional Summary paragraph.
# ---------------------------------------------------------------------
$d.Content.Find.Execute(
    "Alex Morgan is a seasoned Data Scientist with a Bachelor of Science in Computer Science from Stanford University, specializing in machine learning, predictive analytics, and data visualization. With extensive experience in developing advanced machine learning models to enhance business insights and operational efficiency, Alex has a proven track record of leveraging data to drive strategic decision-making and optimize business processes. Proficient in Python, SQL, AWS, and Tableau, Alex is also skilled in statistical analysis and has a strong ability to translate complex datasets into actionable insights. Certified in data science, Alex is eager to contribute to innovative projects at a dynamic company, utilizing technical expertise to solve complex problems and generate significant business value.",
    $true, $false, $false, $false, $false, $true, 1, $false,
    "I am a highly motivated Data Scientist with a Bachelor of Science in Computer Science from Stanford University, specializing in leveraging analytical tools to enhance business outcomes. My expertise encompasses Python, SQL, AWS, and Tableau, with a robust background in machine learning technologies such as Scikit-Learn and TensorFlow. I excel in predictive analytics, including regression analysis and time series forecasting, to improve strategic decision-making processes. With a proven ability to translate complex datasets into actionable business insights, I am adept at driving growth and operational efficiency through data-driven solutions. My professional experience and educational background equip me to contribute effectively to teams focused on innovation and improvement in data science.",
    2) | Out-Null

# ---------------------------------------------------------------------
# 3. Rewrite the Skills paragraph.
# ---------------------------------------------------------------------
$d.Content.Find.Execute(
    "Python, SQL, Machine Learning, Predictive Analytics, Data Visualization, Statistical Analysis, AWS, Tableau, Business Insights",
    $true, $false, $false, $false, $false, $true, 1, $false,
    "Python, SQL, Machine Learning, Predictive Analytics, Data Analysis, AWS, Tableau, Scikit-Learn, TensorFlow, Business Intelligence, Regression Analysis, Time Series Forecasting",
    2) | Out-Null

# ---------------------------------------------------------------------
# 4. Update the "Data Analyst" job heading dates.
# ---------------------------------------------------------------------
$d.Content.Find.Execute(
    "Data Analyst at Tech Solutions (January 2021 - February 2023)",
    $true, $false, $false, $false, $false, $true, 1, $false,
    "Data Analyst at Tech Solutions (Jan 2021 - Feb 2023)",
    2) | Out-Null

# ---------------------------------------------------------------------
# 5. Update the three "Data Analyst" bullet points.
# ---------------------------------------------------------------------
$d.Content.Find.Execute(
    "Developed and maintained advanced dashboards and visualizations using Tableau to track key performance indicators and facilitate strategic decision-making.",
    $true, $false, $false, $false, $false, $true, 1, $false,
    "Developed and maintained advanced dashboards and KPIs to track company metrics, enhancing data visualization and reporting capabilities.",
    2) | Out-Null

$d.Content.Find.Execute(
    "Implemented machine learning models using Python and TensorFlow to analyze large datasets, improving product recommendations and customer segmentation.",
    $true, $false, $false, $false, $false, $true, 1, $false,
    "Performed complex forecasting analysis using statistical tools to predict sales trends, significantly aiding strategic planning and business decisions.",
    2) | Out-Null

$d.Content.Find.Execute(
    "Utilized predictive analytics to forecast market trends and customer behavior, significantly increasing accuracy in targeted marketing campaigns.",
    $true, $false, $false, $false, $false, $true, 1, $false,
    "Implemented machine learning algorithms to optimize and automate data-driven decision-making processes, resulting in a measurable increase in operational efficiency.",
    2) | Out-Null

# ---------------------------------------------------------------------
# 6. Insert a new "Intern Data Scientist" role (Heading2) before the
#    paragraph that used to read "Collaborated with cross-functional...".
#    Then rewrite that paragraph's text to be the first bullet of the
#    new role.
# ---------------------------------------------------------------------
$target = $null
for ($i = 1; $i -le $d.Paragraphs.Count; $i++) {
    $p = $d.Paragraphs($i)
    if ($p.Range.Text.TrimEnd("`r") -eq "Collaborated with cross-functional teams to integrate machine learning algorithms into business processes, enhancing data-driven decision-making across the company.") {
        $target = $p
        break
    }
}
$idx = $target.Index
$target.Range.InsertParagraphBefore()
$newHeading = $d.Paragraphs($idx)
$newHeading.Range.Text = "Intern Data Scientist at Startup Hub (June 2020 - Dec 2020)"
$newHeading.Style = "Heading2"

$d.Content.Find.Execute(
    "Collaborated with cross-functional teams to integrate machine learning algorithms into business processes, enhancing data-driven decision-making across the company.",
    $true, $false, $false, $false, $false, $true, 1, $false,
    "Assisted in the development and refinement of machine learning models to analyze and predict user behavior patterns, increasing model accuracy by 20%.",
    2) | Out-Null

# ---------------------------------------------------------------------
# 7. Rewrite the final bullet and append a brand-new bullet after it.
# ---------------------------------------------------------------------
$d.Content.Find.Execute(
    "Conducted thorough data analysis to identify and interpret trends, providing critical business insights that informed executive-level strategies.",
    $true, $false, $false, $false, $false, $true, 1, $false,
    "Enhanced data collection procedures to include additional relevant information, improving the quality and volume of data for analytics.",
    2) | Out-Null

$lastTarget = $null
for ($i = 1; $i -le $d.Paragraphs.Count; $i++) {
    $p = $d.Paragraphs($i)
    if ($p.Range.Text.TrimEnd("`r") -eq "Enhanced data collection procedures to include additional relevant information, improving the quality and volume of data for analytics.") {
        $lastTarget = $p
        break
    }
}
$lastIdx = $lastTarget.Index
$lastTarget.Range.InsertParagraphAfter()
$newBullet = $d.Paragraphs($lastIdx + 1)
$newBullet.Range.Text = "Collaborated with senior data scientists to create and fine-tune predictive models, boosting user engagement metrics by 15%."
$newBullet.Style = "ListBullet"

Write-Output "done"
